$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = 175136
$ws.Range("C4").Value = 165115
$ws.Range("C5").Value = 10021
$ws.Range("C8").Value = 64.6
